$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated 2D training schedule values (columns B:H for rows 2-6)
# Row 2
$ws.Cells.Item(2, 2).Value = 2
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = -3
$ws.Cells.Item(2, 8).Value = 43

# Row 3
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 5
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = -5
$ws.Cells.Item(3, 8).Value = 21

# Row 4
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = 5
$ws.Cells.Item(4, 7).Value = -1
$ws.Cells.Item(4, 8).Value = 65

# Row 5
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 9
$ws.Cells.Item(5, 4).Value = 5
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(5, 7).Value = -4
$ws.Cells.Item(5, 8).Value = 32

# Row 6
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).Value = 4
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(6, 7).Value = -2
$ws.Cells.Item(6, 8).Value = 54

# Select cell I1, matching the saved selection in the workbook
$ws.Range("I1").Select()
